$d = $word.ActiveDocument

$replacements = @(
    @("38÷3=", "38÷8="),
    @("29÷9=", "16÷9="),
    @("56÷6=", "74÷9="),
    @("66÷7=", "79÷2="),
    @("67÷5=", "84÷9="),
    @("48÷5=", "52÷4="),
    @("48÷4=", "79÷8="),
    @("67÷2=", "42÷2="),
    @("99÷9=", "46÷4="),
    @("90÷7=", "61÷9="),
    @("40÷9=", "34÷8="),
    @("63÷9=", "78÷2="),
    @("93÷3=", "63÷7="),
    @("96÷5=", "97÷5="),
    @("55÷5=", "66÷4="),
    @("19÷2=", "72÷4="),
    @("26÷6=", "33÷6="),
    @("24÷8=", "50÷5="),
    @("44÷5=", "10÷5="),
    @("49÷7=", "50÷2="),
    @("86÷6=", "36÷9="),
    @("96÷6=", "94÷4="),
    @("56÷3=", "85÷5="),
    @("82÷7=", "62÷5="),
    @("17÷6=", "24÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
